$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2..329).
# The source workbook had every row set to serial 46061 (2026-02-08); this
# updates every one of them to 46062 (2026-02-09).
for ($r = 2; $r -le 329; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -eq 46062) {
        continue
    }
    $cell.Value = 46062
}
